# añadido moto2 y moto3 a 2015
# Update Hoja1!B2 from "2014" to "2015" and move the active selection to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B2").Value = "2015"
$ws.Range("B2").Select()
